# Add photo-upload columns (Complaint_ID, Photo_Filename) to the complaints
# sheet and append the new complaint row (row 10 / richard / streetlights).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells G1:H1 -------------------------------------------------
# Copy the header style (bold, bordered, centered - style index 1 in the
# original file) from an existing header cell so the new headers look the
# same as Timestamp/Username/etc.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G1").Value = "Complaint_ID"
$ws.Range("H1").Value = "Photo_Filename"

# --- New (still-blank) G/H cells for the existing data rows (2-9) ----------
# These become "empty" text cells (present, typed as text, but with no
# content) rather than simply unused cells - matching the source data export
# which always writes a Complaint_ID / Photo_Filename cell, blank until a
# photo is attached. Assigning a leading apostrophe forces a text cell with
# empty content, then resetting the style back to Normal drops the
# quote-prefix formatting that the apostrophe would otherwise apply.
foreach ($r in 2..9) {
    $rng = $ws.Range("G" + $r + ":H" + $r)
    $rng.Value = "'"
    $rng.Style = "Normal"
}

# --- New row 10: richard's streetlight complaint, with photo metadata ------
$ws.Range("A10").Value = 45945.78909712595
$ws.Range("B10").Value = "richard"
$ws.Range("C10").Value = "Streetlights in our neighborhood are flickering at night, making it difficult to walk safely. It needs to be checked soon."
$ws.Range("D10").Value = "Roads"
$ws.Range("E10").Value = "Low"
$ws.Range("F10").Value = "at night,, Streetlights, in our"
$ws.Range("G10").Value = "COMP_20251015_185617_5606"
$ws.Range("H10").Value = "complaint_COMP_20251015_185617_5606_20251015_185617_20251015_185617.jpeg"

# Match the date/number formatting of the other rows' Timestamp column.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A10").Value = 45945.78909712595
